# Rebuild the "Equations" matrix with a fresh set of coefficients and an
# extra (5th) row, extending column E so it carries the same number
# format / right alignment as the rest of the table instead of the old
# blank "general" look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new coefficient values (rows 2-5) -------------------------------
$ws.Range("A2").Value = 8
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 24

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 6

$ws.Range("A4").Value = -1
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 1

$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 2

# --- formatting: rows 2,3,5 (cols A-D) take on the header row's look --
$ws.Range("A1:D1").Copy() | Out-Null
$ws.Range("A2:D3").PasteSpecial(-4122) | Out-Null
$ws.Range("A1:D1").Copy() | Out-Null
$ws.Range("A5:D5").PasteSpecial(-4122) | Out-Null

# --- formatting: column E (rows 2-5) now matches the number format /
#     right alignment already used by columns A-D instead of "general"
$ws.Range("A4").Copy() | Out-Null
$ws.Range("E2:E5").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- row heights -------------------------------------------------
$ws.Rows.Item(1).RowHeight = 19.5
$ws.Rows.Item(4).RowHeight = 18.75
$ws.Rows.Item(5).RowHeight = 19.5
